{"js": "// \"first message bug fixed\"\n// Mark the two \"create vehicle / voyage\" TODO sub-items as done:\n//     - no navigation\n//     - isSuccess not defined\n// by applying strikethrough + red font color, matching the formatting\n// already used elsewhere in this document for completed items.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst targets = [\"no navigation\", \"isSuccess not defined\"];\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (targets.some((t) => text.indexOf(t) !== -1)) {\n    paragraph.font.strikeThrough = true;\n    paragraph.font.color = \"#FF0000\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"first message bug fixed\"\n# Mark the two \"create vehicle / voyage\" TODO sub-items as done:\n#     - no navigation\n#     - isSuccess not defined\n# by applying strikethrough + red font color, matching the formatting\n# already used elsewhere in this document for completed items.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -match \"no navigation\" -or $text -match \"isSuccess not defined\") {\n        $r = $p.Range\n        $r.Font.StrikeThrough = 1\n        $r.Font.Color = 255\n    }\n}\n"}
